# Update TPM-derived NATMI ligand-receptor metrics for the Hbegf-Cd9 sheet
# (sending/target cluster expression, specificity and edge-weight columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.459557
$ws.Range("H2").Value = 25.378671
$ws.Range("I2").Value = 0.3030961495696597
$ws.Range("J2").Value = 0.3030961495696597
$ws.Range("M2").Value = 78.90112033333334
$ws.Range("N2").Value = 236.703361
$ws.Range("O2").Value = 0.3371779636489425
$ws.Range("P2").Value = 0.3371779636489425
$ws.Range("Q2").Value = 667.4685248236924
$ws.Range("R2").Value = 6007.216723413231
$ws.Range("S2").Value = 0.1021973425017332
$ws.Range("T2").Value = 0.1021973425017332

# Row 3
$ws.Range("G3").Value = 8.459557
$ws.Range("H3").Value = 25.378671
$ws.Range("I3").Value = 0.3030961495696597
$ws.Range("J3").Value = 0.3030961495696597
$ws.Range("O3").Value = 0.1683260544097508
$ws.Range("P3").Value = 0.1683260544097508
$ws.Range("Q3").Value = 333.213778297938
$ws.Range("R3").Value = 2998.924004681442
$ws.Range("S3").Value = 0.05101897896384851
$ws.Range("T3").Value = 0.0510189789638485

# Row 4
$ws.Range("G4").Value = 8.459557
$ws.Range("H4").Value = 25.378671
$ws.Range("I4").Value = 0.3030961495696597
$ws.Range("J4").Value = 0.3030961495696597
$ws.Range("M4").Value = 44.61912266666666
$ws.Range("N4").Value = 133.857368
$ws.Range("O4").Value = 0.1906764423241422
$ws.Range("P4").Value = 0.1906764423241422
$ws.Range("Q4").Value = 377.4580114886587
$ws.Range("R4").Value = 3397.122103397928
$ws.Range("S4").Value = 0.05779329548208881
$ws.Range("T4").Value = 0.0577932954820888

# Row 5
$ws.Range("G5").Value = 8.459557
$ws.Range("H5").Value = 25.378671
$ws.Range("I5").Value = 0.3030961495696597
$ws.Range("J5").Value = 0.3030961495696597
$ws.Range("M5").Value = 18.49514433333333
$ws.Range("N5").Value = 55.485433
$ws.Range("O5").Value = 0.07903759892585487
$ws.Range("P5").Value = 0.07903759892585487
$ws.Range("Q5").Value = 156.4607277110603
$ws.Range("R5").Value = 1408.146549399543
$ws.Range("S5").Value = 0.02395599190565768
$ws.Range("T5").Value = 0.02395599190565768

# Row 6
$ws.Range("G6").Value = 8.459557
$ws.Range("H6").Value = 25.378671
$ws.Range("I6").Value = 0.3030961495696597
$ws.Range("J6").Value = 0.3030961495696597
$ws.Range("M6").Value = 52.59995866666666
$ws.Range("N6").Value = 157.799876
$ws.Range("O6").Value = 0.2247819406913095
$ws.Range("P6").Value = 0.2247819406913095
$ws.Range("Q6").Value = 444.9723485383106
$ws.Range("R6").Value = 4004.751136844796
$ws.Range("S6").Value = 0.06813054071633154
$ws.Range("T6").Value = 0.06813054071633153

# Row 7
$ws.Range("I7").Value = 0.4601547065605718
$ws.Range("J7").Value = 0.4601547065605718
$ws.Range("M7").Value = 78.90112033333334
$ws.Range("N7").Value = 236.703361
$ws.Range("O7").Value = 0.3371779636489425
$ws.Range("P7").Value = 0.3371779636489425
$ws.Range("Q7").Value = 1013.337792692992
$ws.Range("R7").Value = 9120.040134236926
$ws.Range("S7").Value = 0.1551540269215703
$ws.Range("T7").Value = 0.1551540269215703

# Row 8
$ws.Range("I8").Value = 0.4601547065605718
$ws.Range("J8").Value = 0.4601547065605718
$ws.Range("O8").Value = 0.1683260544097508
$ws.Range("P8").Value = 0.1683260544097508
$ws.Range("S8").Value = 0.07745602617341772
$ws.Range("T8").Value = 0.07745602617341772

# Row 9
$ws.Range("I9").Value = 0.4601547065605718
$ws.Range("J9").Value = 0.4601547065605718
$ws.Range("M9").Value = 44.61912266666666
$ws.Range("N9").Value = 133.857368
$ws.Range("O9").Value = 0.1906764423241422
$ws.Range("P9").Value = 0.1906764423241422
$ws.Range("Q9").Value = 573.0494457356417
$ws.Range("R9").Value = 5157.445011620775
$ws.Range("S9").Value = 0.08774066236567947
$ws.Range("T9").Value = 0.08774066236567946

# Row 10
$ws.Range("I10").Value = 0.4601547065605718
$ws.Range("J10").Value = 0.4601547065605718
$ws.Range("M10").Value = 18.49514433333333
$ws.Range("N10").Value = 55.485433
$ws.Range("O10").Value = 0.07903759892585487
$ws.Range("P10").Value = 0.07903759892585487
$ws.Range("Q10").Value = 237.5356478475812
$ws.Range("R10").Value = 2137.820830628231
$ws.Range("S10").Value = 0.03636952314097891
$ws.Range("T10").Value = 0.03636952314097891

# Row 11
$ws.Range("I11").Value = 0.4601547065605718
$ws.Range("J11").Value = 0.4601547065605718
$ws.Range("M11").Value = 52.59995866666666
$ws.Range("N11").Value = 157.799876
$ws.Range("O11").Value = 0.2247819406913095
$ws.Range("P11").Value = 0.2247819406913095
$ws.Range("Q11").Value = 675.548405217059
$ws.Range("R11").Value = 6079.935646953531
$ws.Range("S11").Value = 0.1034344679589254
$ws.Range("T11").Value = 0.1034344679589254

# Row 12
$ws.Range("G12").Value = 1.955432333333333
$ws.Range("H12").Value = 5.866296999999999
$ws.Range("I12").Value = 0.0700608803720276
$ws.Range("J12").Value = 0.0700608803720276
$ws.Range("M12").Value = 78.90112033333334
$ws.Range("N12").Value = 236.703361
$ws.Range("O12").Value = 0.3371779636489425
$ws.Range("P12").Value = 0.3371779636489425
$ws.Range("Q12").Value = 154.2858018360241
$ws.Range("R12").Value = 1388.572216524217
$ws.Range("S12").Value = 0.02362298497529243
$ws.Range("T12").Value = 0.02362298497529243

# Row 13
$ws.Range("G13").Value = 1.955432333333333
$ws.Range("H13").Value = 5.866296999999999
$ws.Range("I13").Value = 0.0700608803720276
$ws.Range("J13").Value = 0.0700608803720276
$ws.Range("O13").Value = 0.1683260544097508
$ws.Range("P13").Value = 0.1683260544097508
$ws.Range("Q13").Value = 77.022590662366
$ws.Range("R13").Value = 693.2033159612939
$ws.Range("S13").Value = 0.01179307156149696
$ws.Range("T13").Value = 0.01179307156149696

# Row 14
$ws.Range("G14").Value = 1.955432333333333
$ws.Range("H14").Value = 5.866296999999999
$ws.Range("I14").Value = 0.0700608803720276
$ws.Range("J14").Value = 0.0700608803720276
$ws.Range("M14").Value = 44.61912266666666
$ws.Range("N14").Value = 133.857368
$ws.Range("O14").Value = 0.1906764423241422
$ws.Range("P14").Value = 0.1906764423241422
$ws.Range("Q14").Value = 87.2496751473662
$ws.Range("R14").Value = 785.2470763262958
$ws.Range("S14").Value = 0.01335895941543555
$ws.Range("T14").Value = 0.01335895941543555

# Row 15
$ws.Range("G15").Value = 1.955432333333333
$ws.Range("H15").Value = 5.866296999999999
$ws.Range("I15").Value = 0.0700608803720276
$ws.Range("J15").Value = 0.0700608803720276
$ws.Range("M15").Value = 18.49514433333333
$ws.Range("N15").Value = 55.485433
$ws.Range("O15").Value = 0.07903759892585487
$ws.Range("P15").Value = 0.07903759892585487
$ws.Range("Q15").Value = 36.16600323906678
$ws.Range("R15").Value = 325.494029151601
$ws.Range("S15").Value = 0.005537443763236614
$ws.Range("T15").Value = 0.005537443763236614

# Row 16
$ws.Range("G16").Value = 1.955432333333333
$ws.Range("H16").Value = 5.866296999999999
$ws.Range("I16").Value = 0.0700608803720276
$ws.Range("J16").Value = 0.0700608803720276
$ws.Range("M16").Value = 52.59995866666666
$ws.Range("N16").Value = 157.799876
$ws.Range("O16").Value = 0.2247819406913095
$ws.Range("P16").Value = 0.2247819406913095
$ws.Range("Q16").Value = 102.8556599087969
$ws.Range("R16").Value = 925.7009391791718
$ws.Range("S16").Value = 0.01574842065656604
$ws.Range("T16").Value = 0.01574842065656604

# Row 17
$ws.Range("G17").Value = 2.929608
$ws.Range("H17").Value = 8.788824000000002
$ws.Range("I17").Value = 0.1049644685352285
$ws.Range("J17").Value = 0.1049644685352285
$ws.Range("M17").Value = 78.90112033333334
$ws.Range("N17").Value = 236.703361
$ws.Range("O17").Value = 0.3371779636489425
$ws.Range("P17").Value = 0.3371779636489425
$ws.Range("Q17").Value = 231.1493533374961
$ws.Range("R17").Value = 2080.344180037464
$ws.Range("S17").Value = 0.03539170575620184
$ws.Range("T17").Value = 0.03539170575620184

# Row 18
$ws.Range("G18").Value = 2.929608
$ws.Range("H18").Value = 8.788824000000002
$ws.Range("I18").Value = 0.1049644685352285
$ws.Range("J18").Value = 0.1049644685352285
$ws.Range("O18").Value = 0.1683260544097508
$ws.Range("P18").Value = 0.1683260544097508
$ws.Range("Q18").Value = 115.394429118672
$ws.Range("R18").Value = 1038.549862068048
$ws.Range("S18").Value = 0.01766825484175144
$ws.Range("T18").Value = 0.01766825484175144

# Row 19
$ws.Range("G19").Value = 2.929608
$ws.Range("H19").Value = 8.788824000000002
$ws.Range("I19").Value = 0.1049644685352285
$ws.Range("J19").Value = 0.1049644685352285
$ws.Range("M19").Value = 44.61912266666666
$ws.Range("N19").Value = 133.857368
$ws.Range("O19").Value = 0.1906764423241422
$ws.Range("P19").Value = 0.1906764423241422
$ws.Range("Q19").Value = 130.716538717248
$ws.Range("R19").Value = 1176.448848455232
$ws.Range("S19").Value = 0.02001425143074174
$ws.Range("T19").Value = 0.02001425143074173

# Row 20
$ws.Range("G20").Value = 2.929608
$ws.Range("H20").Value = 8.788824000000002
$ws.Range("I20").Value = 0.1049644685352285
$ws.Range("J20").Value = 0.1049644685352285
$ws.Range("M20").Value = 18.49514433333333
$ws.Range("N20").Value = 55.485433
$ws.Range("O20").Value = 0.07903759892585487
$ws.Range("P20").Value = 0.07903759892585487
$ws.Range("Q20").Value = 54.18352280008801
$ws.Range("R20").Value = 487.6517052007921
$ws.Range("S20").Value = 0.008296139565552901
$ws.Range("T20").Value = 0.008296139565552901

# Row 21
$ws.Range("G21").Value = 2.929608
$ws.Range("H21").Value = 8.788824000000002
$ws.Range("I21").Value = 0.1049644685352285
$ws.Range("J21").Value = 0.1049644685352285
$ws.Range("M21").Value = 52.59995866666666
$ws.Range("N21").Value = 157.799876
$ws.Range("O21").Value = 0.2247819406913095
$ws.Range("P21").Value = 0.2247819406913095
$ws.Range("Q21").Value = 154.097259709536
$ws.Range("R21").Value = 1386.875337385824
$ws.Range("S21").Value = 0.02359411694098055
$ws.Range("T21").Value = 0.02359411694098055

# Row 22
$ws.Range("G22").Value = 1.722740333333333
$ws.Range("H22").Value = 5.168221
$ws.Range("I22").Value = 0.06172379496251228
$ws.Range("J22").Value = 0.06172379496251227
$ws.Range("M22").Value = 78.90112033333334
$ws.Range("N22").Value = 236.703361
$ws.Range("O22").Value = 0.3371779636489425
$ws.Range("P22").Value = 0.3371779636489425
$ws.Range("Q22").Value = 135.9261423434201
$ws.Range("R22").Value = 1223.335281090781
$ws.Range("S22").Value = 0.02081190349414475
$ws.Range("T22").Value = 0.02081190349414475

# Row 23
$ws.Range("G23").Value = 1.722740333333333
$ws.Range("H23").Value = 5.168221
$ws.Range("I23").Value = 0.06172379496251228
$ws.Range("J23").Value = 0.06172379496251227
$ws.Range("O23").Value = 0.1683260544097508
$ws.Range("P23").Value = 0.1683260544097508
$ws.Range("Q23").Value = 67.857077562838
$ws.Range("R23").Value = 610.713698065542
$ws.Range("S23").Value = 0.01038972286923614
$ws.Range("T23").Value = 0.01038972286923614

# Row 24
$ws.Range("G24").Value = 1.722740333333333
$ws.Range("H24").Value = 5.168221
$ws.Range("I24").Value = 0.06172379496251228
$ws.Range("J24").Value = 0.06172379496251227
$ws.Range("M24").Value = 44.61912266666666
$ws.Range("N24").Value = 133.857368
$ws.Range("O24").Value = 0.1906764423241422
$ws.Range("P24").Value = 0.1906764423241422
$ws.Range("Q24").Value = 76.86716225581421
$ws.Range("R24").Value = 691.8044603023279
$ws.Range("S24").Value = 0.01176927363019665
$ws.Range("T24").Value = 0.01176927363019665

# Row 25
$ws.Range("G25").Value = 1.722740333333333
$ws.Range("H25").Value = 5.168221
$ws.Range("I25").Value = 0.06172379496251228
$ws.Range("J25").Value = 0.06172379496251227
$ws.Range("M25").Value = 18.49514433333333
$ws.Range("N25").Value = 55.485433
$ws.Range("O25").Value = 0.07903759892585487
$ws.Range("P25").Value = 0.07903759892585487
$ws.Range("Q25").Value = 31.86233111385478
$ws.Range("R25").Value = 286.760980024693
$ws.Range("S25").Value = 0.004878500550428746
$ws.Range("T25").Value = 0.004878500550428746

# Row 26
$ws.Range("G26").Value = 1.722740333333333
$ws.Range("H26").Value = 5.168221
$ws.Range("I26").Value = 0.06172379496251228
$ws.Range("J26").Value = 0.06172379496251227
$ws.Range("M26").Value = 52.59995866666666
$ws.Range("N26").Value = 157.799876
$ws.Range("O26").Value = 0.2247819406913095
$ws.Range("P26").Value = 0.2247819406913095
$ws.Range("Q26").Value = 90.61607032673288
$ws.Range("R26").Value = 815.5446329405959
$ws.Range("S26").Value = 0.01387439441850598
$ws.Range("T26").Value = 0.01387439441850598
